$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 50
$ws.Range("K6").Value = 150
$ws.Range("M6").Value = -38
$ws.Range("H33").Value = 183
$ws.Range("I33").Value = 199.66667
$ws.Range("J33").Value = 133
$ws.Range("K33").Value = 199.66667
$ws.Range("L33").Value = 133
$ws.Range("M33").Value = 29.33332999999999
$ws.Range("N33").Value = -591
$ws.Range("H51").Value = 7176.6924
$ws.Range("J51").Value = 8200
$ws.Range("L51").Value = 8200
$ws.Range("N51").Value = -9168
$ws.Range("H53").Value = 950
$ws.Range("I53").Value = 950
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 950
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -313
$ws.Range("N53").ClearContents()
$ws.Range("H88").Value = 9199.75
$ws.Range("J88").Value = 9199.75
$ws.Range("L88").Value = 9199.75
$ws.Range("N88").Value = -10011.75
$ws.Range("H91").Value = 9199.75
$ws.Range("J91").Value = 9199.75
$ws.Range("L91").Value = 9199.75
$ws.Range("N91").Value = -12007.75
$ws.Range("H100").Value = 2113.5557
$ws.Range("I100").Value = 1943.8235
$ws.Range("K100").Value = 1943.8235
$ws.Range("M100").Value = -1402.8235
$ws.Range("H129").Value = 2628.6365
$ws.Range("I129").Value = 3228
$ws.Range("J129").Value = 2495.4443
$ws.Range("K129").Value = 9684
$ws.Range("L129").Value = 7486.3329
$ws.Range("M129").Value = -4684
$ws.Range("N129").Value = -17486.3329
$ws.Range("H138").Value = 1928.093
$ws.Range("I138").Value = 2083.2222
$ws.Range("J138").Value = 1887.0294
$ws.Range("K138").Value = 6249.6666
$ws.Range("L138").Value = 5661.0882
$ws.Range("M138").Value = -1109.6666
$ws.Range("N138").Value = -15941.0882
$ws.Range("H141").Value = 5243.9
$ws.Range("I141").Value = 6248
$ws.Range("K141").Value = 18744
$ws.Range("M141").Value = -13564

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1908.875
$ws.Range("J2").Value = 1799
$ws.Range("L2").Value = 1799
$ws.Range("N2").Value = -2025
$ws.Range("H32").Value = 1341.0769
$ws.Range("I32").Value = 1341.0769
$ws.Range("K32").Value = 1341.0769
$ws.Range("M32").Value = -1054.0769
$ws.Range("H97").Value = 466.33334
$ws.Range("I97").Value = 466.33334
$ws.Range("K97").Value = 466.33334
$ws.Range("M97").Value = 29.66665999999998
$ws.Range("H102").Value = 2250
$ws.Range("I102").Value = 2250
$ws.Range("K102").Value = 2250
$ws.Range("M102").Value = -628
$ws.Range("H116").Value = 1908.875
$ws.Range("J116").Value = 1799
$ws.Range("L116").Value = 1799
$ws.Range("N116").Value = -6387
$ws.Range("H119").Value = 48999
$ws.Range("J119").Value = 48999
$ws.Range("L119").Value = 48999
$ws.Range("N119").Value = -58675
$ws.Range("H122").Value = 1703.7273
$ws.Range("I122").Value = 1674.1
$ws.Range("K122").Value = 5022.299999999999
$ws.Range("M122").Value = -2572.299999999999
$ws.Range("H132").Value = 4103.75
$ws.Range("I132").Value = 5250.2
$ws.Range("K132").Value = 15750.6
$ws.Range("M132").Value = -13220.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1908.875
$ws.Range("J3").Value = 1799
$ws.Range("L3").Value = 1799
$ws.Range("N3").Value = -2027
$ws.Range("H86").Value = 1584.2
$ws.Range("I86").Value = 1584.2
$ws.Range("K86").Value = 1584.2
$ws.Range("M86").Value = -461.2
$ws.Range("H89").Value = 1584.2
$ws.Range("I89").Value = 1584.2
$ws.Range("K89").Value = 7921
$ws.Range("M89").Value = -2305
$ws.Range("H99").Value = 55558396
$ws.Range("I99").Value = 22223782
$ws.Range("K99").Value = 22223782
$ws.Range("M99").Value = -22222284
$ws.Range("H105").Value = 3752
$ws.Range("I105").Value = 2476.5715
$ws.Range("K105").Value = 2476.5715
$ws.Range("M105").Value = -729.5715
$ws.Range("H134").Value = 1500
$ws.Range("I134").Value = 1500
$ws.Range("K134").Value = 4500
$ws.Range("M134").Value = -1965

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 71430170
$ws.Range("I16").Value = 71430170
$ws.Range("K16").Value = 71430170
$ws.Range("M16").Value = -71429883
$ws.Range("H31").Value = 4009
$ws.Range("I31").Value = 3074.875
$ws.Range("K31").Value = 3074.875
$ws.Range("M31").Value = -2779.875
$ws.Range("H34").Value = 4009
$ws.Range("I34").Value = 3074.875
$ws.Range("K34").Value = 3074.875
$ws.Range("M34").Value = -2872.875
$ws.Range("H58").Value = 4006
$ws.Range("I58").Value = 2012
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 2012
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -1809
$ws.Range("N58").Value = -6406
$ws.Range("H105").Value = 400
$ws.Range("I105").Value = 400
$ws.Range("K105").Value = 400
$ws.Range("M105").Value = 1347
$ws.Range("H113").Value = 71430170
$ws.Range("I113").Value = 71430170
$ws.Range("K113").Value = 71430170
$ws.Range("M113").Value = -71428000
$ws.Range("H122").Value = 1499.5
$ws.Range("J122").Value = 1499.3334
$ws.Range("L122").Value = 4498.0002
$ws.Range("N122").Value = -9398.0002
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
$ws.Range("H134").Value = 3182.0667
$ws.Range("I134").Value = 3240.8462
$ws.Range("K134").Value = 9722.5386
$ws.Range("M134").Value = -7187.5386
$ws.Range("H136").Value = 4006
$ws.Range("I136").Value = 2012
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 6036
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -3486
$ws.Range("N136").Value = -23100

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 6256.5
$ws.Range("I62").Value = 6900
$ws.Range("J62").Value = 5613
$ws.Range("K62").Value = 20700
$ws.Range("L62").Value = 16839
$ws.Range("M62").Value = -20014
$ws.Range("N62").Value = -18211
$ws.Range("H65").Value = 6256.5
$ws.Range("I65").Value = 6900
$ws.Range("J65").Value = 5613
$ws.Range("K65").Value = 62100
$ws.Range("L65").Value = 50517
$ws.Range("M65").Value = -58668
$ws.Range("N65").Value = -57381
$ws.Range("H132").Value = 2059.8
$ws.Range("J132").Value = 2199.75
$ws.Range("L132").Value = 19797.75
$ws.Range("N132").Value = -24857.75
$ws.Range("H133").Value = 17066.625
$ws.Range("I133").Value = 15666.667
$ws.Range("K133").Value = 47000.001
$ws.Range("M133").Value = -41940.001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6179.4
$ws.Range("I5").Value = 6179.4
$ws.Range("K5").Value = 6179.4
$ws.Range("M5").Value = -6067.4
$ws.Range("H126").Value = 2234.4
$ws.Range("I126").Value = 2234.4
$ws.Range("K126").Value = 6703.200000000001
$ws.Range("M126").Value = -4233.200000000001
$ws.Range("H132").Value = 2672.139
$ws.Range("I132").Value = 2583.7097
$ws.Range("K132").Value = 7751.1291
$ws.Range("M132").Value = -5221.1291

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2469.7273
$ws.Range("I46").Value = 1981.2858
$ws.Range("J46").Value = 3324.5
$ws.Range("K46").Value = 1981.2858
$ws.Range("L46").Value = 3324.5
$ws.Range("M46").Value = -1793.2858
$ws.Range("N46").Value = -3700.5
$ws.Range("H136").Value = 35716620
$ws.Range("I136").Value = 2609
$ws.Range("J136").Value = 166668000
$ws.Range("K136").Value = 7827
$ws.Range("L136").Value = 500004000
$ws.Range("M136").Value = -5277
$ws.Range("N136").Value = -500009100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 25027248
$ws.Range("I2").Value = 25027248
$ws.Range("K2").Value = 25027248
$ws.Range("M2").Value = -25027136
$ws.Range("H16").Value = 57710
$ws.Range("J16").Value = 57710
$ws.Range("L16").Value = 57710
$ws.Range("N16").Value = -58294
$ws.Range("H81").Value = 4974.375
$ws.Range("I81").Value = 3959.4
$ws.Range("J81").Value = 6666
$ws.Range("K81").Value = 7918.8
$ws.Range("L81").Value = 13332
$ws.Range("M81").Value = -6857.8
$ws.Range("N81").Value = -15454
$ws.Range("H84").Value = 4974.375
$ws.Range("I84").Value = 3959.4
$ws.Range("J84").Value = 6666
$ws.Range("K84").Value = 39594
$ws.Range("L84").Value = 66660
$ws.Range("M84").Value = -34290
$ws.Range("N84").Value = -77268
$ws.Range("H95").Value = 18368
$ws.Range("J95").Value = 18368
$ws.Range("L95").Value = 18368
$ws.Range("N95").Value = -23860
$ws.Range("H121").Value = 98997
$ws.Range("J121").Value = 98997
$ws.Range("L121").Value = 98997
$ws.Range("N121").Value = -102491
